# Apply diff: restore chronologically-tied match rows to their scraped order,
# and append the newly scraped Fredericia vs Kolding IF match as row 88.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 43: Vendsyssel vs Fredericia
$ws.Cells.Item(43, 6).Value = 'Vendsyssel'
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 'Fredericia'
$ws.Cells.Item(43, 9).Value = 2
$ws.Cells.Item(43, 10).Value = 2.02
$ws.Cells.Item(43, 11).Value = '27/08/2023 13:12'
$ws.Cells.Item(43, 12).Value = 2.46
$ws.Cells.Item(43, 13).Value = '01/09/2023 18:55'
$ws.Cells.Item(43, 14).Value = 3.79
$ws.Cells.Item(43, 15).Value = '27/08/2023 13:12'
$ws.Cells.Item(43, 16).Value = 3.9
$ws.Cells.Item(43, 17).Value = '01/09/2023 18:55'
$ws.Cells.Item(43, 18).Value = 3.19
$ws.Cells.Item(43, 19).Value = '27/08/2023 13:12'
$ws.Cells.Item(43, 20).Value = 2.62
$ws.Cells.Item(43, 21).Value = '01/09/2023 18:55'
$ws.Cells.Item(43, 22).Value = 'https://www.betexplorer.com/football/denmark/1st-division/vendsyssel-ff-fredericia/zmIxmLb8/'

# Row 44: Hillerod vs Sonderjyske
$ws.Cells.Item(44, 6).Value = 'Hillerod'
$ws.Cells.Item(44, 7).Value = 2
$ws.Cells.Item(44, 8).Value = 'Sonderjyske'
$ws.Cells.Item(44, 9).Value = 2
$ws.Cells.Item(44, 10).Value = 3.9
$ws.Cells.Item(44, 11).Value = '28/08/2023 18:42'
$ws.Cells.Item(44, 12).Value = 4.16
$ws.Cells.Item(44, 13).Value = '01/09/2023 18:58'
$ws.Cells.Item(44, 14).Value = 3.8
$ws.Cells.Item(44, 15).Value = '28/08/2023 18:42'
$ws.Cells.Item(44, 16).Value = 3.85
$ws.Cells.Item(44, 17).Value = '01/09/2023 18:58'
$ws.Cells.Item(44, 18).Value = 1.79
$ws.Cells.Item(44, 19).Value = '28/08/2023 18:42'
$ws.Cells.Item(44, 20).Value = 1.83
$ws.Cells.Item(44, 21).Value = '01/09/2023 18:58'
$ws.Cells.Item(44, 22).Value = 'https://www.betexplorer.com/football/denmark/1st-division/hillerod-sonderjyske/EwHtnuEE/'

# Row 45: Horsens vs B.93
$ws.Cells.Item(45, 6).Value = 'Horsens'
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = 'B.93'
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 1.79
$ws.Cells.Item(45, 11).Value = '28/08/2023 01:42'
$ws.Cells.Item(45, 12).Value = 1.59
$ws.Cells.Item(45, 13).Value = '01/09/2023 18:57'
$ws.Cells.Item(45, 14).Value = 4.02
$ws.Cells.Item(45, 15).Value = '28/08/2023 01:42'
$ws.Cells.Item(45, 16).Value = 4.54
$ws.Cells.Item(45, 17).Value = '01/09/2023 18:57'
$ws.Cells.Item(45, 18).Value = 3.72
$ws.Cells.Item(45, 19).Value = '28/08/2023 01:42'
$ws.Cells.Item(45, 20).Value = 4.97
$ws.Cells.Item(45, 21).Value = '01/09/2023 18:57'
$ws.Cells.Item(45, 22).Value = 'https://www.betexplorer.com/football/denmark/1st-division/horsens-boldklubben-1893/MqUfOyM7/'

# Row 61: Vendsyssel vs Aalborg
$ws.Cells.Item(61, 6).Value = 'Vendsyssel'
$ws.Cells.Item(61, 7).Value = 1
$ws.Cells.Item(61, 8).Value = 'Aalborg'
$ws.Cells.Item(61, 9).Value = 3
$ws.Cells.Item(61, 10).Value = 3.87
$ws.Cells.Item(61, 11).Value = '23/09/2023 17:13'
$ws.Cells.Item(61, 12).Value = 4.85
$ws.Cells.Item(61, 13).Value = '29/09/2023 18:46'
$ws.Cells.Item(61, 14).Value = 3.84
$ws.Cells.Item(61, 15).Value = '23/09/2023 17:13'
$ws.Cells.Item(61, 16).Value = 4.09
$ws.Cells.Item(61, 17).Value = '29/09/2023 18:47'
$ws.Cells.Item(61, 18).Value = 1.85
$ws.Cells.Item(61, 19).Value = '23/09/2023 17:13'
$ws.Cells.Item(61, 20).Value = 1.67
$ws.Cells.Item(61, 21).Value = '29/09/2023 18:45'
$ws.Cells.Item(61, 22).Value = 'https://www.betexplorer.com/football/denmark/1st-division/vendsyssel-ff-aalborg/0E7GucVE/'

# Row 62: Hillerod vs Horsens
$ws.Cells.Item(62, 6).Value = 'Hillerod'
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 8).Value = 'Horsens'
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 2.44
$ws.Cells.Item(62, 11).Value = '24/09/2023 16:13'
$ws.Cells.Item(62, 12).Value = 2.98
$ws.Cells.Item(62, 13).Value = '29/09/2023 18:32'
$ws.Cells.Item(62, 14).Value = 3.52
$ws.Cells.Item(62, 15).Value = '24/09/2023 16:13'
$ws.Cells.Item(62, 16).Value = 3.55
$ws.Cells.Item(62, 17).Value = '29/09/2023 18:32'
$ws.Cells.Item(62, 18).Value = 2.66
$ws.Cells.Item(62, 19).Value = '24/09/2023 16:13'
$ws.Cells.Item(62, 20).Value = 2.34
$ws.Cells.Item(62, 21).Value = '29/09/2023 18:32'
$ws.Cells.Item(62, 22).Value = 'https://www.betexplorer.com/football/denmark/1st-division/hillerod-horsens/jqvSYQd2/'

# Row 80: Hobro vs Koge
$ws.Cells.Item(80, 6).Value = 'Hobro'
$ws.Cells.Item(80, 7).Value = 2
$ws.Cells.Item(80, 8).Value = 'Koge'
$ws.Cells.Item(80, 9).Value = 1
$ws.Cells.Item(80, 10).Value = 1.74
$ws.Cells.Item(80, 11).Value = '22/10/2023 16:12'
$ws.Cells.Item(80, 12).Value = 1.65
$ws.Cells.Item(80, 13).Value = '27/10/2023 18:51'
$ws.Cells.Item(80, 14).Value = 3.92
$ws.Cells.Item(80, 15).Value = '22/10/2023 16:12'
$ws.Cells.Item(80, 16).Value = 4.26
$ws.Cells.Item(80, 17).Value = '27/10/2023 18:51'
$ws.Cells.Item(80, 18).Value = 4.42
$ws.Cells.Item(80, 19).Value = '22/10/2023 16:12'
$ws.Cells.Item(80, 20).Value = 4.85
$ws.Cells.Item(80, 21).Value = '27/10/2023 18:51'
$ws.Cells.Item(80, 22).Value = 'https://www.betexplorer.com/football/denmark/1st-division/hobro-koge/UBGxoh7f/'

# Row 81: B.93 vs Sonderjyske
$ws.Cells.Item(81, 6).Value = 'B.93'
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = 'Sonderjyske'
$ws.Cells.Item(81, 9).Value = 4
$ws.Cells.Item(81, 10).Value = 5.01
$ws.Cells.Item(81, 11).Value = '22/10/2023 15:12'
$ws.Cells.Item(81, 12).Value = 8.029999999999999
$ws.Cells.Item(81, 13).Value = '27/10/2023 18:58'
$ws.Cells.Item(81, 14).Value = 4.47
$ws.Cells.Item(81, 15).Value = '22/10/2023 15:12'
$ws.Cells.Item(81, 16).Value = 5.46
$ws.Cells.Item(81, 17).Value = '27/10/2023 18:58'
$ws.Cells.Item(81, 18).Value = 1.57
$ws.Cells.Item(81, 19).Value = '22/10/2023 15:12'
$ws.Cells.Item(81, 20).Value = 1.34
$ws.Cells.Item(81, 21).Value = '27/10/2023 18:50'
$ws.Cells.Item(81, 22).Value = 'https://www.betexplorer.com/football/denmark/1st-division/boldklubben-1893-sonderjyske/2PEpqWy7/'

# New row 88: Fredericia vs Kolding IF (newly scraped match)
$ws.Range("A87").Copy()
$ws.Range("A88").PasteSpecial(-4122)
$ws.Range("E87").Copy()
$ws.Range("E88").PasteSpecial(-4122)

$ws.Cells.Item(88, 1).Value = 87
$ws.Cells.Item(88, 2).Value = 'denmark'
$ws.Cells.Item(88, 3).Value = '1st-division'
$ws.Cells.Item(88, 4).Value = '2023-2024'
$ws.Cells.Item(88, 5).Value = 45235.58333333334
$ws.Cells.Item(88, 6).Value = 'Fredericia'
$ws.Cells.Item(88, 7).Value = 1
$ws.Cells.Item(88, 8).Value = 'Kolding IF'
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 2.15
$ws.Cells.Item(88, 11).Value = '29/10/2023 14:12'
$ws.Cells.Item(88, 12).Value = 2.59
$ws.Cells.Item(88, 13).Value = '05/11/2023 13:53'
$ws.Cells.Item(88, 14).Value = 3.68
$ws.Cells.Item(88, 15).Value = '29/10/2023 14:12'
$ws.Cells.Item(88, 16).Value = 3.64
$ws.Cells.Item(88, 17).Value = '05/11/2023 13:53'
$ws.Cells.Item(88, 18).Value = 3.15
$ws.Cells.Item(88, 19).Value = '29/10/2023 14:12'
$ws.Cells.Item(88, 20).Value = 2.61
$ws.Cells.Item(88, 21).Value = '05/11/2023 13:53'
$ws.Cells.Item(88, 22).Value = 'https://www.betexplorer.com/football/denmark/1st-division/fredericia-kolding-if/M3O4lRDn/'
